$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("collectLevel")

# Update the "dev of alliance fight" cost column (D) to a flat 600 for every level row
$ws.Range("D2").Value = 600
$ws.Range("D3").Value = 600
$ws.Range("D4").Value = 600
$ws.Range("D5").Value = 600
$ws.Range("D6").Value = 600
$ws.Range("D7").Value = 600
$ws.Range("D8").Value = 600
$ws.Range("D9").Value = 600
$ws.Range("D10").Value = 600
$ws.Range("D11").Value = 600

# Move the active selection to E9 (was E10)
$ws.Activate()
$ws.Range("E9").Select()
